$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, centered alignment) that the
# existing header cell H1 already carries over to the two new header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats = -4122
$excel.CutCopyMode = 0

# --- Data rows 2-5: add I column (constant 1) and J column (mirrors H) ---
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = $ws.Range("H2").Value2

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = $ws.Range("H3").Value2

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = $ws.Range("H4").Value2

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = $ws.Range("H5").Value2
